# Cotações atualizadas - 2025-11-29
# Adds a new row (85) with the quotes for 2025-11-29 (Excel serial date 45990),
# mirroring the layout of the existing rows (date in column A with the same
# date number format, values as text in columns B:E using comma decimals).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 85
$lastRow = $newRow - 1

# Column A: numeric date serial, formatted like the row above it.
$ws.Cells.Item($newRow, 1).Value = 45990
$ws.Cells.Item($newRow, 1).NumberFormat = $ws.Cells.Item($lastRow, 1).NumberFormat

# Columns B:E: textual values (comma decimal separator), as in the source data.
$ws.Cells.Item($newRow, 2).Value = "21,7883"
$ws.Cells.Item($newRow, 3).Value = "16,0515"
$ws.Cells.Item($newRow, 4).Value = "15,5122"
$ws.Cells.Item($newRow, 5).Value = "15,5122"
